$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: "Save" - match the style of the other header cells (G1 etc.)
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Column H values per row (0/1 flag), matching the diff
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
